# Apply the diff: insert one new row right after the header (row 2),
# shifting all existing data rows down by one, and insert another new
# row just before the final existing row (which, after the first
# insertion, sits at row 11), shifting that last row down to row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextRow($rowIndex, $phone, $ddd, $date) {
    # Write values with a leading apostrophe so Excel treats them as text
    # (preventing "+55..." from being parsed as a number and the date
    # string from being parsed as a date serial).
    $ws.Range("A" + $rowIndex).Value = "'" + $phone
    $ws.Range("B" + $rowIndex).Value = "'" + $ddd
    $ws.Range("C" + $rowIndex).Value = "'" + $date

    # Re-apply the same formatting used by the rest of the data rows by
    # copying the format from the row directly below (an existing,
    # correctly-styled data row) onto the new row. This is done after
    # setting the values so the text values themselves are preserved.
    $ws.Range("A" + ($rowIndex + 1) + ":C" + ($rowIndex + 1)).Copy()
    $ws.Range("A" + $rowIndex + ":C" + $rowIndex).PasteSpecial(-4122)
}

# 1) Insert the new first data row (was not present before) at row 2.
$ws.Rows.Item(2).Insert()
Set-TextRow 2 "+553171858800" "31" "2024-09-30"

# 2) Insert another new row just before the last existing row, which is
#    now located at row 11 (old row 10 shifted down by the insert above).
$ws.Rows.Item(11).Insert()
Set-TextRow 11 "+5521985109311" "21" "2024-09-09"

$ws.Range("A1").Select()
